# Update the two-digit multiplication answer table with a new set of
# generated problems/answers. Each data row of the 5-column table holds
# one "A×B=C" string per cell; only the w:t text changes, formatting is
# left untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row (1-based, Word Table.Cell indexing) -> list of 5 new cell values,
# in column order, matching the order of replacements in the diff.
$replacements = @{
    1  = @("31×59=1829", "53×48=2544", "92×52=4784", "65×21=1365", "33×69=2277")
    5  = @("70×92=6440", "33×36=1188", "79×46=3634", "37×76=2812", "50×81=4050")
    10 = @("35×81=2835", "28×94=2632", "22×64=1408", "31×60=1860", "50×77=3850")
    15 = @("95×57=5415", "92×33=3036", "12×16=192",  "80×12=960",  "12×26=312")
    20 = @("79×20=1580", "48×19=912",  "78×39=3042", "85×98=8330", "42×65=2730")
}

foreach ($row in $replacements.Keys) {
    $values = $replacements[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cellRange = $t.Cell($row, $col).Range
        $cellRange.MoveEnd(1, -1) | Out-Null
        $cellRange.Text = $values[$col - 1]
    }
}
